$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.335.57'
$ws.Range('E2').Value = '  -1.41%  '
$ws.Range('D3').Value = '2.577.26'
$ws.Range('E3').Value = '  -2.78%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'588.62"
$ws.Range('E5').Value = '  -3.26%  '
$ws.Range('D6').Value = "'150.40"
$ws.Range('E6').Value = '  +0.90%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  -0.64%  '
$ws.Range('E9').Value = '  +1.08%  '
$ws.Range('D10').Value = "'5.71"
$ws.Range('E10').Value = '  +1.93%  '
$ws.Range('E11').Value = '  -0.78%  '
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('D14').Value = '3.038.44'
$ws.Range('E14').Value = '  -2.63%  '
$ws.Range('D15').Value = '63.173.31'
$ws.Range('D16').Value = "'0.0000155"
$ws.Range('E16').Value = '  +4.97%  '
$ws.Range('D17').Value = '2.593.11'
$ws.Range('E17').Value = '  -2.97%  '
$ws.Range('E18').Value = '  +2.22%  '
$ws.Range('E19').Value = '  +3.03%  '
$ws.Range('D20').Value = "'345.61"
$ws.Range('E20').Value = '  -0.32%  '
$ws.Range('D21').Value = "'6.84"
$ws.Range('E21').Value = '  -1.12%  '
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').Value = "'67.08"
$ws.Range('E23').Value = '  +1.15%  '
$ws.Range('E24').Value = '  +1.35%  '
$ws.Range('D25').Value = "'1.67"
$ws.Range('E25').Value = '  -3.76%  '
$ws.Range('D26').Value = "'9.12"
$ws.Range('E26').Value = '  -2.91%  '
$ws.Range('D27').Value = "'553.35"
$ws.Range('E27').Value = '  -0.29%  '
$ws.Range('E28').Value = '  -2.06%  '
$ws.Range('E29').Value = '  +0.52%  '
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('D31').Value = "'2.04"
$ws.Range('E31').Value = '  -1.31%  '
$ws.Range('D32').Value = '0.0₃0855'
$ws.Range('E32').Value = '  +0.22%  '
$ws.Range('D33').Value = "'1.75"
$ws.Range('E33').Value = '  -0.53%  '
$ws.Range('E34').Value = '  -1.55%  '
$ws.Range('D35').Value = "'166.52"
$ws.Range('E35').Value = '  -2.00%  '
$ws.Range('D36').Value = "'0.412"
$ws.Range('E36').Value = '  +1.15%  '
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('E38').Value = '  +0.72%  '
$ws.Range('E39').Value = '  -1.91%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').Value = "'165.18"
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('D42').Value = "'39.67"
$ws.Range('E42').Value = '  -1.53%  '
$ws.Range('E43').Value = '  +3.43%  '
$ws.Range('E44').Value = '  +3.41%  '
$ws.Range('D45').Value = "'0.0584"
$ws.Range('E45').Value = '  +2.34%  '
$ws.Range('E46').Value = '  +4.99%  '
$ws.Range('E47').Value = '  -0.31%  '
$ws.Range('E48').Value = '  +1.70%  '
$ws.Range('D49').Value = "'0.0961"
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('D50').Value = "'19.15"
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('D51').Value = '0.0₆0234'
$ws.Range('E51').Value = '  +19.37%  '
